$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "26.948.41"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "1.817.69"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  +0.18%  "
Set-TextValue $ws.Range("D5") "310.26"
$ws.Range("E5").Value = "  +0.24%  "
Set-TextValue $ws.Range("D6") "1.003"
$ws.Range("E6").Value = "  +0.16%  "
Set-TextValue $ws.Range("D7") "0.4686"
$ws.Range("E7").Value = "  +0.95%  "
Set-TextValue $ws.Range("D8") "0.3665"
$ws.Range("E8").Value = "  -0.80%  "
Set-TextValue $ws.Range("D9") "0.07345"
$ws.Range("E9").Value = "  -0.09%  "
Set-TextValue $ws.Range("D10") "0.8723"
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("D12").Value = "1.822.61"
$ws.Range("E12").Value = "  -1.70%  "
Set-TextValue $ws.Range("D13") "5.404"
$ws.Range("E13").Value = "  +0.98%  "
Set-TextValue $ws.Range("D14") "0.07112"
Set-TextValue $ws.Range("D15") "6.511"
$ws.Range("E15").Value = "  -0.06%  "
Set-TextValue $ws.Range("D16") "91.35"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("D21").Value = "26.966.79"
$ws.Range("E21").Value = "  +0.23%  "
Set-TextValue $ws.Range("D22") "5.281"
$ws.Range("E22").Value = "  -0.70%  "
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("D24").Value = "2.047.22"
$ws.Range("E24").Value = "  -1.38%  "
Set-TextValue $ws.Range("D25") "1.893"
$ws.Range("E25").Value = "  -0.56%  "
Set-TextValue $ws.Range("D26") "150.94"
Set-TextValue $ws.Range("D27") "18.36"
$ws.Range("E27").Value = "  -0.03%  "
Set-TextValue $ws.Range("D28") "2.152"
$ws.Range("E28").Value = "  +0.52%  "
Set-TextValue $ws.Range("D29") "5.244"
$ws.Range("E29").Value = "  -1.11%  "
$ws.Range("E30").Value = "  +1.07%  "
Set-TextValue $ws.Range("D31") "0.08891"
$ws.Range("E31").Value = "  -0.03%  "
Set-TextValue $ws.Range("D32") "0.7589"
$ws.Range("E32").Value = "  +0.71%  "
$ws.Range("E33").Value = "  +0.73%  "
Set-TextValue $ws.Range("D34") "4.505"
$ws.Range("E34").Value = "  +1.12%  "
Set-TextValue $ws.Range("D35") "2.913"
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("E37").Value = "  +0.06%  "
Set-TextValue $ws.Range("D38") "0.05309"
$ws.Range("E38").Value = "  +1.10%  "
$ws.Range("E39").Value = "  -0.74%  "
Set-TextValue $ws.Range("D40") "2.966"
$ws.Range("E40").Value = "  +1.36%  "
Set-TextValue $ws.Range("D41") "2.378"
$ws.Range("E41").Value = "  -1.43%  "
Set-TextValue $ws.Range("D42") "0.5294"
$ws.Range("E42").Value = "  -0.29%  "
Set-TextValue $ws.Range("D43") "7.166"
$ws.Range("E43").Value = "  +0.08%  "
Set-TextValue $ws.Range("D44") "0.1655"
$ws.Range("E44").Value = "  -0.22%  "
Set-TextValue $ws.Range("D45") "8.436"
$ws.Range("E45").Value = "  +0.03%  "
Set-TextValue $ws.Range("D46") "0.4871"
$ws.Range("E46").Value = "  -1.36%  "
Set-TextValue $ws.Range("D47") "10.49"
$ws.Range("E47").Value = "  +2.12%  "
$ws.Range("E48").Value = "  +0.15%  "
Set-TextValue $ws.Range("D49") "103.38"
$ws.Range("E49").Value = "  +0.29%  "
Set-TextValue $ws.Range("D50") "1.661"
$ws.Range("E50").Value = "  -0.43%  "
Set-TextValue $ws.Range("D51") "0.06301"
$ws.Range("E51").Value = "  +0.32%  "
